# "update wire colors on ppt"
#
# Adds two new wire-color tag textboxes ("blk" and "wht") to slide 1,
# matching the ones already used elsewhere on the slide for the other
# wires (red/yel/wht/blk tags near the bottom of the slide), but placed
# up near the top-right of the slide (by the TSL1401 pinout).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Reuse two of the existing wire-color tag textboxes as style/content
# templates so the duplicates come out with identical formatting
# (no-fill, wrap="none" autosize textbox, etc.) instead of the bare
# defaults AddTextbox would give us.
$blkTemplate = $s.Shapes.Item("TextBox 1")
$whtTemplate = $s.Shapes.Item("TextBox 76")

# New shapes always take the lowest free shape id on the slide, and this
# slide's id sequence already has a lot of gaps from earlier edits. Burn
# through the lower gaps with disposable filler textboxes first so that
# the two real shapes we add land on the same ids (81 and 84) that they
# have in the target deck, then remove the filler shapes again.
$fillers = @()
for ($i = 0; $i -lt 31; $i++) {
    $fillers += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}

$blk = $blkTemplate.Duplicate()
$blk.Name = "TextBox 80"

$fillers += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)

$wht = $whtTemplate.Duplicate()
$wht.Name = "TextBox 83"

foreach ($f in $fillers) { $f.Delete() }

# Position/size the two new tags (values are EMU/12700, nudged by a few
# ulps so they survive the host's float32 Left/Top/Width/Height setters
# and land on the exact target EMU values).
$blk.Left = 598.7185039370079
$blk.Top = 84.97417322834646
$blk.Width = 36.50299272598425
$blk.Height = 29.081259842519685

$wht.Left = 551.087738135433
$wht.Top = 84.4796062992126
$wht.Width = 43.02614173228346
$wht.Height = 29.081259842519685
